$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pl_mw results table (rows 2-25) for columns B, C, D, F, G, H, K, M
# with the new simulation results for the 380 kV case.
$data = @{
    "2" = @{ "B"=0.2162994289648736; "C"=0.03133084960579424; "D"=0.03393069819610872; "F"=1.138836170362822; "G"=1.010047872608311; "H"=0.9857809319671844; "K"=0.1708619479876887; "M"=0.8801035284143666 }
    "3" = @{ "B"=0.1945429446079174; "C"=0.02965730065874084; "D"=0.03360386571623764; "F"=1.116958623902804; "G"=0.987655369128575; "H"=0.9792795701159491; "K"=0.1500835450877531; "M"=0.7839656560886397 }
    "4" = @{ "B"=0.1812701132699459; "C"=0.02861051222124189; "D"=0.03339213200879421; "F"=1.104061161900148; "G"=0.9743881710405589; "H"=0.9756552959098883; "K"=0.137349619966102; "M"=0.7254832238757984 }
    "5" = @{ "B"=0.1758831402260057; "C"=0.02817912639009279; "D"=0.03330308439064211; "F"=1.098939502468468; "G"=0.9691022251439705; "H"=0.9742706183521221; "K"=0.1321667547295675; "M"=0.7017821167133178 }
    "6" = @{ "B"=0.1749899624732194; "C"=0.02810720519035925; "D"=0.03328813160917576; "F"=1.098097143875677; "G"=0.9682317601336052; "H"=0.9740462594888868; "K"=0.1313065326671961; "M"=0.6978542880536338 }
    "7" = @{ "B"=0.1811973740068424; "C"=0.02860471385039176; "D"=0.03339094225341377; "F"=1.103991546824773; "G"=0.9743163956703569; "H"=0.9756362484189083; "K"=0.137279696136531; "M"=0.7251630607030251 }
    "8" = @{ "B"=0.2087801427575187; "C"=0.03075781275224898; "D"=0.03382031081264003; "F"=1.131181363466482; "G"=1.002226532768759; "H"=0.9834628512678449; "K"=0.1636926977161295; "M"=0.8468377719657525 }
    "9" = @{ "B"=0.2635419861717878; "C"=0.03482666994077732; "D"=0.03457387825182145; "F"=1.188774508677923; "G"=1.060813933605317; "H"=1.001737679011995; "K"=0.2156712960070877; "M"=1.090087323697361 }
    "10" = @{ "B"=0.3041786633279742; "C"=0.03772166002668342; "D"=0.03507271384838972; "F"=1.233734694984136; "G"=1.106258496999999; "H"=1.01696512470852; "K"=0.2539640291128933; "M"=1.272085978004071 }
    "11" = @{ "B"=0.3227517873469594; "C"=0.03901798711891757; "D"=0.03528755548003204; "F"=1.254772141603283; "G"=1.127464762513171; "H"=1.024287314088554; "K"=0.2714056822293287; "M"=1.355696024925507 }
    "12" = @{ "B"=0.3297973129845957; "C"=0.03950588659493803; "D"=0.03536715696783332; "F"=1.262823212364509; "G"=1.135572510131567; "H"=1.027117121457422; "K"=0.2780133644091336; "M"=1.387482755126499 }
    "13" = @{ "B"=0.3282793907455641; "C"=0.03940094218302193; "D"=0.03535009166124325; "F"=1.261085496761552; "G"=1.133822908124728; "H"=1.026505130631364; "K"=0.2765901568472486; "M"=1.380631209561642 }
    "14" = @{ "B"=0.3233311822569362; "C"=0.03905818700837926; "D"=0.03529413960046668; "F"=1.255432808441824; "G"=1.128130236606381; "H"=1.02451897927736; "K"=0.2719492435116422; "M"=1.358308579326533 }
    "15" = @{ "B"=0.3203018538963249; "C"=0.03884784919906537; "D"=0.03525963837003943; "F"=1.2519814134439; "G"=1.124653411134972; "H"=1.023309841152098; "K"=0.2691069206113923; "M"=1.344651892175094 }
    "16" = @{ "B"=0.3029665997913753; "C"=0.03763652483961977; "D"=0.03505842887318167; "F"=1.232371663412579; "G"=1.104883409383092; "H"=1.016494573501006; "K"=0.2528245975245511; "M"=1.266638986116178 }
    "17" = @{ "B"=0.2923541512597012; "C"=0.03688811685873361; "D"=0.0349318868751034; "F"=1.220491894641967; "G"=1.0928922465651; "H"=1.012414997992749; "K"=0.2428413860669423; "M"=1.218995576601401 }
    "18" = @{ "B"=0.2862583840328909; "C"=0.03645571327342623; "D"=0.03485796714437228; "F"=1.213713979934425; "G"=1.086045401510205; "H"=1.010105713239682; "K"=0.2371014028268377; "M"=1.191668857761115 }
    "19" = @{ "B"=0.2841958840761549; "C"=0.03630897661304999; "D"=0.03483274458682928; "F"=1.211428526585237; "G"=1.083735770105932; "H"=1.009330207658024; "K"=0.2351583113840263; "M"=1.182429406673748 }
    "20" = @{ "B"=0.2934830144273519; "C"=0.03696798704744708; "D"=0.0349454751956273; "F"=1.221750818609493; "G"=1.094163530407798; "H"=1.012845426202119; "K"=0.2439039010129704; "M"=1.224059324508445 }
    "21" = @{ "B"=0.3247842596930468; "C"=0.03915894391099073; "D"=0.0353106218264152; "F"=1.257090838349384; "G"=1.129800206079352; "H"=1.02510080980764; "K"=0.2733123155205988; "M"=1.364861809887472 }
    "22" = @{ "B"=0.3453129267905126; "C"=0.04057341192175556; "D"=0.03553903252669599; "F"=1.280681226362205; "G"=1.153542284341427; "H"=1.033443070730527; "K"=0.2925492320622425; "M"=1.457619186023237 }
    "23" = @{ "B"=0.3343499366134211; "C"=0.03982008912041834; "D"=0.03541806745735343; "F"=1.268045241307092; "G"=1.140829147397596; "H"=1.028960133395827; "K"=0.2822806779791165; "M"=1.408042962014918 }
    "24" = @{ "B"=0.2929726383285924; "C"=0.03693188437938488; "D"=0.03493933555608564; "F"=1.22118149739768; "G"=1.093588636535031; "H"=1.012650717128338; "K"=0.2434235395489992; "M"=1.221769804972794 }
    "25" = @{ "B"=0.248656191963164; "C"=0.03374244819838879; "D"=0.03437957863046748; "F"=1.172732260225644; "G"=1.044546575488681; "H"=0.9964788027899658; "K"=0.2015909090158061; "M"=1.023738538602785 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
